$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Pixel clock changed from a computed 113.75/2 (56.875 MHz) to a fixed 57 MHz value.
$ws.Range("B1").Value = 57

# Update selection to match the author's saved cursor position.
$ws.Range("F24").Select()
